$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$ws.Range("D2").Value = "26.338.47"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").Value = "1.666.17"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  +0.97%  "

Set-TextValue $ws.Range("D5") "219.40"
$ws.Range("E5").Value = "  +1.00%  "

Set-TextValue $ws.Range("D6") "0.5343"
$ws.Range("E6").Value = "  +1.58%  "

$ws.Range("E7").Value = "  +0.89%  "

Set-TextValue $ws.Range("D8") "0.2660"

Set-TextValue $ws.Range("D9") "0.06384"
$ws.Range("E9").Value = "  +1.26%  "

Set-TextValue $ws.Range("D10") "20.83"
$ws.Range("E10").Value = "  +2.65%  "

Set-TextValue $ws.Range("D11") "0.07838"
$ws.Range("E11").Value = "  +0.59%  "

Set-TextValue $ws.Range("D12") "4.557"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").Value = "1.672.38"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").Value = "1.894.52"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "0.0₅8179"
$ws.Range("E16").Value = "  +0.08%  "

Set-TextValue $ws.Range("D17") "65.90"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "26.359.61"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("E19").Value = "  +0.90%  "

Set-TextValue $ws.Range("D20") "4.671"
$ws.Range("E20").Value = "  +2.26%  "

Set-TextValue $ws.Range("D21") "193.79"
$ws.Range("E21").Value = "  +1.87%  "

Set-TextValue $ws.Range("D22") "10.27"
$ws.Range("E22").Value = "  +1.94%  "

Set-TextValue $ws.Range("D23") "6.036"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("E24").Value = "  +0.94%  "

Set-TextValue $ws.Range("D25") "145.91"
$ws.Range("E25").Value = "  +1.89%  "

$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("E27").Value = "  -0.16%  "

Set-TextValue $ws.Range("D28") "16.09"
$ws.Range("E28").Value = "  +0.68%  "

Set-TextValue $ws.Range("D29") "1.500"
$ws.Range("E29").Value = "  +4.64%  "

Set-TextValue $ws.Range("D30") "0.05857"
$ws.Range("E30").Value = "  +0.75%  "

Set-TextValue $ws.Range("D31") "1.282"
$ws.Range("E31").Value = "  +1.05%  "

Set-TextValue $ws.Range("D32") "3.593"
$ws.Range("E32").Value = "  +1.39%  "

Set-TextValue $ws.Range("D33") "3.283"
$ws.Range("E33").Value = "  +0.82%  "

Set-TextValue $ws.Range("D34") "1.601"
$ws.Range("E34").Value = "  +0.87%  "

Set-TextValue $ws.Range("D35") "0.9678"
$ws.Range("E35").Value = "  +2.81%  "

$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("E37").Value = "  +0.33%  "

Set-TextValue $ws.Range("D38") "0.5811"
$ws.Range("E38").Value = "  +1.29%  "

Set-TextValue $ws.Range("D39") "0.01608"
$ws.Range("E39").Value = "  +0.31%  "

Set-TextValue $ws.Range("D40") "0.8622"
$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("D41").Value = "1.065.98"
$ws.Range("E41").Value = "  +3.64%  "

Set-TextValue $ws.Range("D42") "5.832"
$ws.Range("E42").Value = "  +2.26%  "

$ws.Range("E43").Value = "  +0.94%  "

Set-TextValue $ws.Range("D44") "104.21"
$ws.Range("E44").Value = "  -0.98%  "

$ws.Range("D45").Value = "1.804.89"
$ws.Range("E45").Value = "  +0.72%  "

Set-TextValue $ws.Range("D46") "57.93"
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D47") "1.017"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.4394"
$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "8.046"
$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈102"
$ws.Range("E50").Value = "  -8.35%  "

Set-TextValue $ws.Range("D51") "0.05167"
$ws.Range("E51").Value = "  +0.53%  "
